# "Generate Report for Handoff"
#
# The localization status report moves from "In Translation" to
# "Ready for handoff", with fresh handoff timestamps, on all three sheets
# (Overview, zh-cn, de-de). The status/date columns also widen slightly to
# fit the new "Ready for handoff" label.
#
# Note on the column width: the target OOXML column width is
# 17.2159881591797 "characters". Excel's Range.ColumnWidth is quantized to
# a whole-pixel grid (MDW=7 -> 1 px = 1/6 character) before it is stored,
# so the nearest value actually reachable through the ColumnWidth property
# is 17.1666... (round(17.2159881591797 * 6) / 6). Any input in
# [16.25, 16.4166...) collapses to that same stored width, so 16.33 is used
# as a safely-centered value.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-01 21:07:57"
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-09-01 21:07:51"
$wsZh.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-09-01 21:07:57"
$wsDe.Columns.Item(3).ColumnWidth = 16.33
